$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "About Me" title paragraph so the new paragraph can be
# anchored reliably (rather than assuming a fixed paragraph index).
# ------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq "About Me") {
        $titlePara = $p
        break
    }
}

if ($titlePara -eq $null) {
    throw "Could not find the 'About Me' title paragraph"
}

$followingPara = $titlePara.Next()

# Insert a brand-new, blank paragraph right before the paragraph that
# currently follows the title - this lands the new paragraph directly
# under "About Me" without inheriting the title's language formatting.
$followingPara.Range.InsertParagraphBefore()

# The freshly inserted blank paragraph is now titlePara's immediate
# successor.
$newPara = $titlePara.Next()

# ------------------------------------------------------------------
# Fill the new paragraph with the two runs exactly as authored: the
# text was typed as "I am a web dev" followed by a second insertion
# completing "eloper who enjoys ...", so it is represented as two
# separate <w:r> runs in the canonical OOXML.
# ------------------------------------------------------------------
$run1 = "I am a web dev"
$run2 = "eloper who enjoys using the tools of the web to craft impressive interactive sites and applications that satisfy the users needs. I began my journey as a self-taught developer, and from then I keep feeding myself with resources to hone my coding skills and keep a modern outlook on the web. I am constantly evolving, and I enjoy working on artistic, colorful projects that bring joy to people, as well as practical programming that increases efficiency and automates repetitive tasks.  "

$packageXml = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>$run1</w:t></w:r>
<w:r><w:t xml:space="preserve">$run2</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$null = $newPara.Range.InsertXML($packageXml)
